$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark (bookmarkStart/bookmarkEnd) if present.
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2. Update the mailing address.
$d.Content.Find.Execute(
    "3014 Whispering Willow Way, Orange Park, Fl, 32065",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2516 Ambrosia Drive, Middleburg, Florida, 32068",
    2)

# 3. Update the email line: split "Email: edavis0314@hotmail.com" so the
#    address is in its own run, with the new address "ecdavis0314@gmail.com".
$d.Content.Find.Execute(
    "edavis0314@hotmail.com",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "ecdavis0314@gmail.com",
    2)

# Force "Email: " to become its own run (separate from the address run) by
# toggling a formatting property on and back off over just that text.
$prefix = $d.Content
$prefix.Find.Execute("Email: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$prefix.Bold = 1
$prefix.Bold = 0

$d.Save()
